$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price / Volume(1h) columns store numeric-looking values as literal text
# (e.g. "35.40", "2.99%") so trailing zeros and the "%" sign survive exactly.
# Force each touched cell to Text format before writing so Excel does not
# silently convert the string into a number/percentage.

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2.99%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "35.40"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-1.97%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.138"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.71%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08128"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "3.19%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.136"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.03%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4.148"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.76%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "7.961"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.14%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9302"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.21%"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "3.75%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1872"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.57%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09066"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "3.66%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03604"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "1.42%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09910"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.02%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001438"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.85%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005683"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.33%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.465"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.08%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.762"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "4.91%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3412"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.54%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1340"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.93%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.097"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.41%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "9.97%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04559"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.17%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001247"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.92%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004706"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-6.64%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-21.92%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-5.13%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01959"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "5.86%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04860"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.08%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007712"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "2.35%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.65%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "1.18%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002172"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-1.12%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01182"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "6.84%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006610"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "4.52%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.05%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "37.71"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-21.66%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-14.88%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.05%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.05%"

# Plain text cells (coin name / link) - rows 7 and 8 swapped their contents.
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
